$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "On shadow map, compute light projection based on actual camera position
#  or give an option to keep it fixed"
# splits into two runs: the first part gets a green highlight, the rest
# (" or give an option to keep it fixed") keeps its original formatting.
$rng1 = $d.Content
$rng1.Find.Execute("On shadow map, compute light projection based on actual camera position", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng1.Find.Found) {
    $rng1.Font.HighlightColorIndex = 4   # wdBrightGreen -> w:highlight val="green"
}

# --- Change 2 -----------------------------------------------------------
# The "BUG" (red highlight) run and the ": shadow only material partially
# works on Android" run merge into a single run, and the whole thing
# (including the paragraph mark) becomes green-highlighted instead of red.

# Remove the second half of the text, then re-append it right after "BUG"
# so it becomes part of that same run (keeping its rPr, e.g. w:lang).
$rngTail = $d.Content
$rngTail.Find.Execute(": shadow only material partially works on Android", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rngTail.Find.Found) {
    $rngTail.Delete()
}

$rngBug = $d.Content
$rngBug.Find.Execute("BUG", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rngBug.Find.Found) {
    $rngBug.InsertAfter(": shadow only material partially works on Android")
}

# Re-highlight the merged paragraph (text + paragraph mark) green.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "BUG: shadow only material partially works on Android*") {
        $p.Range.Font.HighlightColorIndex = 4   # wdBrightGreen -> w:highlight val="green"
    }
}
